# Applies cryptos list price/volume update (GitHub Actions scheduled refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($cellRange, [string]$text)
    # Leading apostrophe forces Excel to store the value as literal text
    # instead of auto-coercing numeric-looking strings (e.g. "0.672",
    # "0.1000") into numbers. Resetting the style back to Normal afterwards
    # strips the quote-prefix formatting Excel stamps on forced-text cells,
    # so the cell keeps the original (unstyled) look.
    $cellRange.Value = "'" + $text
    $cellRange.Style = "Normal"
}

Set-TextValue $ws.Range('D2') '43.741.01'
Set-TextValue $ws.Range('E2') '  -0.89%  '
Set-TextValue $ws.Range('D3') '2.345.93'
Set-TextValue $ws.Range('E3') '  -0.33%  '
Set-TextValue $ws.Range('E4') '  +0.17%  '
Set-TextValue $ws.Range('D5') '0.672'
Set-TextValue $ws.Range('E5') '  -0.92%  '
Set-TextValue $ws.Range('D6') '239.15'
Set-TextValue $ws.Range('E6') '  -0.84%  '
Set-TextValue $ws.Range('D7') '73.12'
Set-TextValue $ws.Range('E7') '  -1.25%  '
Set-TextValue $ws.Range('E8') '  +0.00%  '
Set-TextValue $ws.Range('D9') '0.598'
Set-TextValue $ws.Range('E9') '  +6.04%  '
Set-TextValue $ws.Range('D10') '0.1000'
Set-TextValue $ws.Range('E10') '  -2.21%  '
Set-TextValue $ws.Range('D11') '58.60'
Set-TextValue $ws.Range('E11') '  +2.54%  '
Set-TextValue $ws.Range('D12') '32.70'
Set-TextValue $ws.Range('E12') '  +4.19%  '
Set-TextValue $ws.Range('D13') '7.32'
Set-TextValue $ws.Range('E13') '  -1.22%  '
Set-TextValue $ws.Range('D14') '0.107'
Set-TextValue $ws.Range('E14') '  -0.35%  '
Set-TextValue $ws.Range('D15') '2.695.58'
Set-TextValue $ws.Range('E15') '  -0.29%  '
Set-TextValue $ws.Range('D16') '16.34'
Set-TextValue $ws.Range('E16') '  -2.91%  '
Set-TextValue $ws.Range('D17') '0.900'
Set-TextValue $ws.Range('E17') '  -1.11%  '
Set-TextValue $ws.Range('D18') '2.343.17'
Set-TextValue $ws.Range('E18') '  -0.53%  '
Set-TextValue $ws.Range('D19') '43.660.64'
Set-TextValue $ws.Range('E19') '  -1.71%  '
Set-TextValue $ws.Range('E20') '  -1.27%  '
Set-TextValue $ws.Range('D21') '6.71'
Set-TextValue $ws.Range('E21') '  +0.29%  '
Set-TextValue $ws.Range('D22') '77.29'
Set-TextValue $ws.Range('E22') '  -0.82%  '
Set-TextValue $ws.Range('D23') '253.65'
Set-TextValue $ws.Range('E23') '  -0.91%  '
Set-TextValue $ws.Range('D24') '1.95'
Set-TextValue $ws.Range('E24') '  +22.67%  '
Set-TextValue $ws.Range('E25') '  -0.05%  '
Set-TextValue $ws.Range('D26') '3.74'
Set-TextValue $ws.Range('E26') '  -0.89%  '
Set-TextValue $ws.Range('E27') '  -2.91%  '
Set-TextValue $ws.Range('D28') '10.61'
Set-TextValue $ws.Range('E28') '  -0.30%  '
Set-TextValue $ws.Range('B29') 'Toncoin'
Set-TextValue $ws.Range('C29') 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
Set-TextValue $ws.Range('D29') '2.27'
Set-TextValue $ws.Range('E29') '  -1.89%  '
Set-TextValue $ws.Range('B30') 'EthereumClassic'
Set-TextValue $ws.Range('C30') 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
Set-TextValue $ws.Range('D30') '22.62'
Set-TextValue $ws.Range('E30') '  +0.29%  '
Set-TextValue $ws.Range('D31') '177.05'
Set-TextValue $ws.Range('E31') '  +1.44%  '
Set-TextValue $ws.Range('D32') '0.130'
Set-TextValue $ws.Range('E32') '  -1.10%  '
Set-TextValue $ws.Range('E33') '  +3.18%  '
Set-TextValue $ws.Range('D34') '0.0758'
Set-TextValue $ws.Range('E34') '  +1.11%  '
Set-TextValue $ws.Range('E35') '  -3.69%  '
Set-TextValue $ws.Range('D36') '5.48'
Set-TextValue $ws.Range('E36') '  +2.47%  '
Set-TextValue $ws.Range('D37') '3.82'
Set-TextValue $ws.Range('E37') '  -2.11%  '
Set-TextValue $ws.Range('D38') '2.37'
Set-TextValue $ws.Range('E38') '  -3.43%  '
Set-TextValue $ws.Range('D39') '6.27'
Set-TextValue $ws.Range('E39') '  -4.12%  '
Set-TextValue $ws.Range('D40') '0.0282'
Set-TextValue $ws.Range('E40') '  +2.98%  '
Set-TextValue $ws.Range('D41') '68.57'
Set-TextValue $ws.Range('E41') '  +30.37%  '
Set-TextValue $ws.Range('E42') '  +11.53%  '
Set-TextValue $ws.Range('B43') 'FraxShare'
Set-TextValue $ws.Range('C43') 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
Set-TextValue $ws.Range('D43') '9.13'
Set-TextValue $ws.Range('E43') '  +1.59%  '
Set-TextValue $ws.Range('B44') 'Algorand'
Set-TextValue $ws.Range('C44') 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
Set-TextValue $ws.Range('D44') '0.203'
Set-TextValue $ws.Range('E44') '  +7.12%  '
Set-TextValue $ws.Range('D45') '19.02'
Set-TextValue $ws.Range('E45') '  -2.14%  '
Set-TextValue $ws.Range('D46') '4.85'
Set-TextValue $ws.Range('E46') '  +8.13%  '
Set-TextValue $ws.Range('D47') '2.51'
Set-TextValue $ws.Range('E47') '  +2.28%  '
Set-TextValue $ws.Range('E48') '  +0.23%  '
Set-TextValue $ws.Range('E49') '  -1.96%  '
Set-TextValue $ws.Range('D50') '99.28'
Set-TextValue $ws.Range('E50') '  -0.92%  '
Set-TextValue $ws.Range('E51') '  -0.83%  '
